$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 350, shifting all existing rows (350..430)
# down to (352..432). This matches the weekly data refresh: a new week's
# worth of "Albahaca" (Primera/Segunda) observations is inserted at the
# top of this date-ordered block.
$ws.Rows.Item(350).Insert()
$ws.Rows.Item(350).Insert()

# New row 350: Primera
$ws.Cells.Item(350, 1).Value = 6
$ws.Cells.Item(350, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(350, 3).Value = "Metropolitana"
$ws.Cells.Item(350, 4).Value = 44637
$ws.Cells.Item(350, 5).Value = 13
$ws.Cells.Item(350, 6).Value = 100112052
$ws.Cells.Item(350, 7).Value = "Albahaca"
$ws.Cells.Item(350, 8).Value = "Sin especificar"
$ws.Cells.Item(350, 9).Value = "Primera"
$ws.Cells.Item(350, 10).Value = 620
$ws.Cells.Item(350, 11).Value = 2500
$ws.Cells.Item(350, 12).Value = 3000
$ws.Cells.Item(350, 13).Value = 2815
$ws.Cells.Item(350, 14).Value = "`$/docena de matas"
$ws.Cells.Item(350, 15).Value = "Región Metropolitana"
$ws.Cells.Item(350, 16).Value = 469
$ws.Cells.Item(350, 17).Value = 6
$ws.Cells.Item(350, 18).Value = "Hortaliza"

# New row 351: Segunda
$ws.Cells.Item(351, 1).Value = 6
$ws.Cells.Item(351, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(351, 3).Value = "Metropolitana"
$ws.Cells.Item(351, 4).Value = 44637
$ws.Cells.Item(351, 5).Value = 13
$ws.Cells.Item(351, 6).Value = 100112052
$ws.Cells.Item(351, 7).Value = "Albahaca"
$ws.Cells.Item(351, 8).Value = "Sin especificar"
$ws.Cells.Item(351, 9).Value = "Segunda"
$ws.Cells.Item(351, 10).Value = 260
$ws.Cells.Item(351, 11).Value = 2000
$ws.Cells.Item(351, 12).Value = 2500
$ws.Cells.Item(351, 13).Value = 2308
$ws.Cells.Item(351, 14).Value = "`$/docena de matas"
$ws.Cells.Item(351, 15).Value = "Región Metropolitana"
$ws.Cells.Item(351, 16).Value = 385
$ws.Cells.Item(351, 17).Value = 6
$ws.Cells.Item(351, 18).Value = "Hortaliza"
